# Commit: "Changed org role to org affiliation"
#
# The diagram slide has two rounded-rectangle boxes whose label is split
# across two paragraphs: "Organization" / "Role" and "Practitioner" /
# "Role". Only the Organization box's second line changes, from "Role"
# to "Affiliation" (OrganizationRole -> OrganizationAffiliation); the
# Practitioner/Role box is left as-is.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Organization" / "Role" shape by its text rather than a
# hard-coded index, in case shape ordering ever shifts.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        $candidateText = $candidate.TextFrame.TextRange.Text
        if ($candidateText -eq "OrganizationRole") {
            $target = $candidate
            break
        }
    }
}

if ($target -eq $null) {
    # Fallback: "Rounded Rectangle 11" (shape id 12) is known to hold it.
    $target = $s.Shapes.Item(2)
}

$tr = $target.TextFrame.TextRange

# The text box has two paragraphs ("Organization" then "Role"/"Affiliation").
# Replace just the second paragraph's characters so the existing run
# formatting (size, color, etc.) carries over unchanged.
$secondPara = $tr.Paragraphs(2, 1)
$secondRun = $tr.Characters($secondPara.Start, $secondPara.Length)
$secondRun.Text = "Affiliation"
